$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$newSlide = $s.Duplicate()
$newSlide.Item(1).MoveTo(2)
